$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; existing rows 16..77 shift down to 17..78.
$ws.Range("A16").EntireRow.Insert()

# Populate the newly inserted row 16 with the new data record.
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Vega Modelo de Temuco"
$ws.Range("C16").Value = "La Araucanía"
$ws.Range("D16").Value = 44462
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 100112012
$ws.Range("G16").Value = "Espinaca"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 12000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = 12000
$ws.Range("N16").Value = "$/docena de atados"
$ws.Range("O16").Value = "Región de La Araucanía"
$ws.Range("P16").Value = 4000
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = "Hortaliza"
